$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header in A1 from "khach_hang" to "maKH"
$ws.Range("A1").Value = "maKH"

# Move selection to A2 to match the saved state
$ws.Range("A2").Select()
